$d = $word.ActiveDocument

# Replace the date line and each table-cell arithmetic expression with new values.
# Using Find/Execute with MatchWholeWord to be precise and MatchCase to avoid ambiguity.
$replacements = @(
    ,@("2024-03-15 Friday", "2024-03-16 Saturday")
    ,@("25+35=60", "91-52=39")
    ,@("54-37=17", "32+43=75")
    ,@("79-2=77", "12+6=18")
    ,@("43+27=70", "97-63=34")
    ,@("30+63=93", "59+10=69")
    ,@("65-27=38", "76-69=7")
    ,@("60-19=41", "76-23=53")
    ,@("28+45=73", "40+26=66")
    ,@("3+47=50", "30-10=20")
    ,@("90-50=40", "89-57=32")
    ,@("58-37=21", "77+12=89")
    ,@("82-53=29", "55+15=70")
    ,@("28+9=37", "75-65=10")
    ,@("37+31=68", "73-5=68")
    ,@("25+41=66", "17-4=13")
    ,@("16+62=78", "93-31=62")
    ,@("83-52=31", "49+20=69")
    ,@("84-1=83", "80-17=63")
    ,@("21+67=88", "5+88=93")
    ,@("1+7=8", "47-44=3")
    ,@("91-46=45", "30-23=7")
    ,@("66-49=17", "38-20=18")
    ,@("76-45=31", "28+31=59")
    ,@("54+23=77", "56-25=31")
    ,@("57+12=69", "68-43=25")
    ,@("43+14=57", "33-18=15")
    ,@("61-32=29", "67+14=81")
    ,@("65+34=99", "14+63=77")
    ,@("2+45=47", "56+18=74")
    ,@("82-40=42", "37+60=97")
    ,@("10+85=95", "32-0=32")
    ,@("12+11=23", "95-17=78")
    ,@("61-33=28", "69-44=25")
    ,@("23-12=11", "31+34=65")
    ,@("26+57=83", "38-9=29")
    ,@("92-36=56", "28+41=69")
    ,@("35-4=31", "31+48=79")
    ,@("13+69=82", "5+93=98")
    ,@("19+61=80", "83-61=22")
    ,@("29+36=65", "62+28=90")
    ,@("36+37=73", "33+61=94")
    ,@("90-54=36", "94-45=49")
    ,@("63+7=70", "19-17=2")
    ,@("48-46=2", "45-24=21")
    ,@("32-24=8", "79-21=58")
    ,@("10-7=3", "38+35=73")
    ,@("12+49=61", "17+76=93")
    ,@("10+33=43", "1+6=7")
    ,@("94-78=16", "52+36=88")
    ,@("64-48=16", "82-0=82")
    ,@("60-59=1", "78+19=97")
    ,@("94-58=36", "99-41=58")
    ,@("74+23=97", "25+17=42")
    ,@("99-82=17", "8+47=55")
    ,@("58-51=7", "37+0=37")
    ,@("18+8=26", "12-5=7")
    ,@("11+31=42", "49-34=15")
    ,@("97-68=29", "28+57=85")
    ,@("77-33=44", "92-0=92")
    ,@("20+12=32", "29+24=53")
    ,@("29+65=94", "78-21=57")
    ,@("7+80=87", "93-65=28")
    ,@("31+40=71", "45+34=79")
    ,@("70+7=77", "39-3=36")
    ,@("95-71=24", "85-1=84")
    ,@("28+23=51", "33+12=45")
    ,@("41+13=54", "18+27=45")
    ,@("47-43=4", "0+23=23")
    ,@("22+13=35", "13+8=21")
    ,@("13+42=55", "52-48=4")
    ,@("40+38=78", "36+23=59")
    ,@("34+28=62", "21-13=8")
    ,@("45-30=15", "24-17=7")
    ,@("22-6=16", "67+10=77")
    ,@("50+24=74", "30+59=89")
    ,@("34+16=50", "55-46=9")
    ,@("59-32=27", "28+39=67")
    ,@("46+52=98", "8+33=41")
    ,@("26-23=3", "66+11=77")
    ,@("75-58=17", "7-3=4")
    ,@("72-68=4", "5+81=86")
    ,@("54-33=21", "44+25=69")
    ,@("37+34=71", "87-85=2")
    ,@("4+40=44", "68+1=69")
    ,@("96-23=73", "9+3=12")
    ,@("67-48=19", "0+56=56")
    ,@("91-35=56", "26+35=61")
    ,@("82-14=68", "51+24=75")
    ,@("59+16=75", "51-47=4")
    ,@("36-8=28", "5+76=81")
    ,@("21+34=55", "68-65=3")
    ,@("96-58=38", "97-88=9")
    ,@("20+76=96", "48+5=53")
    ,@("8+48=56", "14+67=81")
    ,@("2+73=75", "94-24=70")
    ,@("70-64=6", "68-20=48")
    ,@("22+41=63", "22+43=65")
    ,@("96-46=50", "4+20=24")
    ,@("50-32=18", "89-43=46")
    ,@("79-50=29", "74-49=25")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
